$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 50052.5
$ws.Range("I8").Value = 105
$ws.Range("J8").Value = 100000
$ws.Range("K8").Value = 315
$ws.Range("L8").Value = 300000
$ws.Range("M8").Value = -176
$ws.Range("N8").Value = -300278
$ws.Range("H112").Value = 1802.8108
$ws.Range("J112").Value = 1903.2122
$ws.Range("L112").Value = 5709.6366
$ws.Range("N112").Value = -7925.6366
$ws.Range("H137").Value = 10601.98
$ws.Range("I137").Value = 5153.727
$ws.Range("J137").Value = 14882.75
$ws.Range("K137").Value = 15461.181
$ws.Range("L137").Value = 44648.25
$ws.Range("M137").Value = -12911.181
$ws.Range("N137").Value = -49748.25
$ws.Range("H138").Value = 3800.756
$ws.Range("I138").Value = 5549
$ws.Range("J138").Value = 3440.8235
$ws.Range("K138").Value = 16647
$ws.Range("L138").Value = 10322.4705
$ws.Range("M138").Value = -11507
$ws.Range("N138").Value = -20602.4705
$ws.Range("H141").Value = 2964.8125
$ws.Range("I141").Value = 2962.4666
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 8887.399800000001
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -3707.399800000001
$ws.Range("N141").Value = -19360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 72502500
$ws.Range("I11").Value = 96668340
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 96668340
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -96668196
$ws.Range("N11").Value = -5288
$ws.Range("H32").Value = 6668.7065
$ws.Range("I32").Value = 5701.409
$ws.Range("J32").Value = 27949.25
$ws.Range("K32").Value = 5701.409
$ws.Range("L32").Value = 27949.25
$ws.Range("M32").Value = -5414.409
$ws.Range("N32").Value = -28523.25
$ws.Range("H37").Value = 32216.889
$ws.Range("J37").Value = 32216.889
$ws.Range("L37").Value = 32216.889
$ws.Range("N37").Value = -32762.889
$ws.Range("H61").Value = 20749
$ws.Range("I61").Value = 13665
$ws.Range("J61").Value = 27833
$ws.Range("K61").Value = 13665
$ws.Range("L61").Value = 27833
$ws.Range("M61").Value = -13453
$ws.Range("N61").Value = -28257
$ws.Range("H88").Value = 1824.6
$ws.Range("J88").Value = 2062.0625
$ws.Range("L88").Value = 2062.0625
$ws.Range("N88").Value = -2874.0625
$ws.Range("H91").Value = 1824.6
$ws.Range("J91").Value = 2062.0625
$ws.Range("L91").Value = 2062.0625
$ws.Range("N91").Value = -4870.0625
$ws.Range("H136").Value = 20749
$ws.Range("I136").Value = 13665
$ws.Range("J136").Value = 27833
$ws.Range("K136").Value = 40995
$ws.Range("L136").Value = 83499
$ws.Range("M136").Value = -38445
$ws.Range("N136").Value = -88599

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 68749.75
$ws.Range("J58").Value = 65000
$ws.Range("L58").Value = 65000
$ws.Range("N58").Value = -65588
$ws.Range("H60").Value = 65000
$ws.Range("J60").Value = 65000
$ws.Range("L60").Value = 65000
$ws.Range("N60").Value = -66198
$ws.Range("H86").Value = 336789.4
$ws.Range("I86").Value = 910361.8
$ws.Range("K86").Value = 910361.8
$ws.Range("M86").Value = -909238.8
$ws.Range("H89").Value = 336789.4
$ws.Range("I89").Value = 910361.8
$ws.Range("K89").Value = 4551809
$ws.Range("M89").Value = -4546193

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2924.4375
$ws.Range("I31").Value = 2345.25
$ws.Range("J31").Value = 4662
$ws.Range("K31").Value = 2345.25
$ws.Range("L31").Value = 4662
$ws.Range("M31").Value = -2050.25
$ws.Range("N31").Value = -5252
$ws.Range("H34").Value = 2924.4375
$ws.Range("I34").Value = 2345.25
$ws.Range("J34").Value = 4662
$ws.Range("K34").Value = 2345.25
$ws.Range("L34").Value = 4662
$ws.Range("M34").Value = -2143.25
$ws.Range("N34").Value = -5066
$ws.Range("H58").Value = 4497.3184
$ws.Range("I58").Value = 2734.606
$ws.Range("J58").Value = 9785.454
$ws.Range("K58").Value = 2734.606
$ws.Range("L58").Value = 9785.454
$ws.Range("M58").Value = -2531.606
$ws.Range("N58").Value = -10191.454
$ws.Range("H88").Value = 26665
$ws.Range("J88").Value = 26665
$ws.Range("L88").Value = 26665
$ws.Range("N88").Value = -27477
$ws.Range("H91").Value = 26665
$ws.Range("J91").Value = 26665
$ws.Range("L91").Value = 26665
$ws.Range("N91").Value = -29473
$ws.Range("H132").Value = 21542.621
$ws.Range("I132").Value = 15077.255
$ws.Range("K132").Value = 45231.765
$ws.Range("M132").Value = -42701.765
$ws.Range("H136").Value = 4497.3184
$ws.Range("I136").Value = 2734.606
$ws.Range("J136").Value = 9785.454
$ws.Range("K136").Value = 8203.818000000001
$ws.Range("L136").Value = 29356.362
$ws.Range("M136").Value = -5653.818000000001
$ws.Range("N136").Value = -34456.362

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22412954
$ws.Range("I4").Value = 23577452
$ws.Range("J4").Value = 5333666.5
$ws.Range("K4").Value = 70732356
$ws.Range("L4").Value = 16000999.5
$ws.Range("M4").Value = -70732244
$ws.Range("N4").Value = -16001223.5
$ws.Range("H7").Value = 166666750
$ws.Range("I7").Value = 166666750
$ws.Range("K7").Value = 500000250
$ws.Range("M7").Value = -500000138
$ws.Range("H92").Value = 3112.25
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H131").Value = 2857.1428
$ws.Range("I131").Value = 1339.1111
$ws.Range("J131").Value = 3271.1516
$ws.Range("K131").Value = 4017.3333
$ws.Range("L131").Value = 9813.4548
$ws.Range("M131").Value = 1022.6667
$ws.Range("N131").Value = -19893.4548

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 25596.477
$ws.Range("J43").Value = 34253.285
$ws.Range("L43").Value = 34253.285
$ws.Range("N43").Value = -34555.285
$ws.Range("H57").Value = 36918.92
$ws.Range("J57").Value = 36918.92
$ws.Range("L57").Value = 36918.92
$ws.Range("N57").Value = -38558.92
$ws.Range("H80").Value = 2958.4
$ws.Range("J80").Value = 2995
$ws.Range("L80").Value = 2995
$ws.Range("N80").Value = -4991
$ws.Range("H83").Value = 2958.4
$ws.Range("J83").Value = 2995
$ws.Range("L83").Value = 14975
$ws.Range("N83").Value = -24959
$ws.Range("H134").Value = 38155.332
$ws.Range("J134").Value = 38155.332
$ws.Range("L134").Value = 114465.996
$ws.Range("N134").Value = -119535.996

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 3356
$ws.Range("I8").Value = 3820
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 3820
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -3680
$ws.Range("N8").Value = -1780
$ws.Range("H45").Value = 10411.857
$ws.Range("J45").Value = 13949.5
$ws.Range("L45").Value = 13949.5
$ws.Range("N45").Value = -14931.5
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H96").Value = 2422
$ws.Range("I96").Value = 1740.0625
$ws.Range("J96").Value = 5149.75
$ws.Range("K96").Value = 1740.0625
$ws.Range("L96").Value = 5149.75
$ws.Range("M96").Value = -367.0625
$ws.Range("N96").Value = -7895.75
$ws.Range("H107").Value = 2105.8572
$ws.Range("I107").Value = 1407.7273
$ws.Range("J107").Value = 4665.6665
$ws.Range("K107").Value = 4223.1819
$ws.Range("L107").Value = 13996.9995
$ws.Range("M107").Value = -2303.1819
$ws.Range("N107").Value = -17836.9995
$ws.Range("H122").Value = 5047.6665
$ws.Range("I122").Value = 5019.25
$ws.Range("K122").Value = 15057.75
$ws.Range("M122").Value = -12607.75
$ws.Range("H132").Value = 126530.29
$ws.Range("I132").Value = 162785.31
$ws.Range("J132").Value = 30817
$ws.Range("K132").Value = 488355.93
$ws.Range("L132").Value = 92451
$ws.Range("M132").Value = -485825.93
$ws.Range("N132").Value = -97511

Write-Host "All updates applied"
